$d = $word.ActiveDocument

function FindParaByPrefix($searchText) {
    foreach ($p in $d.Paragraphs) {
        if ($p.Range.Text.StartsWith($searchText)) {
            return $p
        }
    }
    return $null
}

function SplitRunAt($absPos, $paraEnd) {
    # Force a w:r run boundary exactly at $absPos (without altering visible formatting)
    # by toggling Bold on/off for the range from $absPos to the end of the paragraph text.
    $r = $d.Range($absPos, $paraEnd)
    $r.Font.Bold = 1
    $r.Font.Bold = 0
}

# ---------------------------------------------------------------------------
# Step 1: Insert a brand new bullet paragraph for "V 1.02" before the existing
#         "v 1.01 - 30 September 2012 - Bug Fix Release." bullet. The new
#         paragraph inherits the list/paragraph formatting automatically
#         because InsertBefore at the start of that paragraph's range picks up
#         the paragraph mark/format of the following paragraph.
# ---------------------------------------------------------------------------
$pOld101 = FindParaByPrefix("v 1.01 - 30 September 2012")
$insertRange = $d.Range($pOld101.Range.Start, $pOld101.Range.Start)
$insertRange.InsertBefore("V 1.02 " + [char]0x2013 + " 8 October 2012 " + [char]0x2013 + " Bug Fix Release. " + [char]13)

# ---------------------------------------------------------------------------
# Step 2: Rewrite the (now shifted) "v 1.01 ..." paragraph so that its text
#         reads "v 1.01 - 30 September 2012 - Bug Fix Release." ->
#         "v 1.01 <en dash> 30 September 2012 <en dash> Bug Fix Release."
#         and split it into the 5 target runs: "v" / " 1.01 <en dash> " /
#         "30" / " September 2012" / " <en dash> Bug Fix Release."
# ---------------------------------------------------------------------------
$p101 = FindParaByPrefix("v 1.01 - 30 September 2012")
$full101 = $p101.Range
$start101 = $full101.Start
$textRange101 = $d.Range($full101.Start, $full101.End - 1)
$textRange101.Text = "v 1.01 " + [char]0x2013 + " 30 September 2012 " + [char]0x2013 + " Bug Fix Release."

# Re-resolve the paragraph end after the text replacement (length changed).
$p101b = FindParaByPrefix("v 1.01")
$paraEnd101 = $p101b.Range.End - 1

foreach ($offset in @(26, 11, 9, 1)) {
    SplitRunAt ($start101 + $offset) $paraEnd101
}

# ---------------------------------------------------------------------------
# Step 3: Rewrite the leading portion of the (now shifted) "v 1.0 - 28 ..."
#         paragraph (up through "Derived from Mike "), leaving the rest of
#         the paragraph (Magatagan's / hyperlink / trailing period) untouched.
#         New text: "v 1.0  <en dash> 28 September 2012  <en dash>" +
#                    " Initial Release. Derived from Mike "
# ---------------------------------------------------------------------------
$p10 = FindParaByPrefix("v 1.0 - 28 September 2012")
$full10 = $p10.Range
$sub10 = $d.Range($full10.Start, $full10.Start)
$sub10.End = $full10.End
$sub10.Find.Execute("v 1.0 - 28 September 2012 - Initial Release. Derived from Mike ") | Out-Null
$start10 = $sub10.Start

$sub10.Text = "v 1.0  " + [char]0x2013 + " 28 September 2012  " + [char]0x2013 + " Initial Release. Derived from Mike "

$p10b = FindParaByPrefix("v 1.0  ")
$paraEnd10 = $p10b.Range.End - 1

SplitRunAt ($start10 + 29) $paraEnd10

$d.Save()
